$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Trim trailing whitespace from the ProductName values for the
# "Streamliner Roadster" and "Carbonite Century" rows.
$ws.Range("B3").Value = "Streamliner Roadster"
$ws.Range("B4").Value = "Carbonite Century"

# Update the selected cell/range to match the saved selection (B9).
$ws.Range("B9").Select()
